$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("I2").Value = 0.08983953209358128
$ws.Range("J2").Value = 0.0898395320935813
$ws.Range("M2").Value = 3.867218333333334
$ws.Range("N2").Value = 11.601655
$ws.Range("O2").Value = 0.1566152977872902
$ws.Range("P2").Value = 0.1566152977872902
$ws.Range("Q2").Value = 0.04633314285166667
$ws.Range("R2").Value = 0.4169982856650001
$ws.Range("S2").Value = 0.01407024507190705
$ws.Range("T2").Value = 0.01407024507190705

$ws.Range("I3").Value = 0.08983953209358128
$ws.Range("J3").Value = 0.0898395320935813
$ws.Range("N3").Value = 33.813685
$ws.Range("O3").Value = 0.4564642152831324
$ws.Range("P3").Value = 0.4564642152831324
$ws.Range("S3").Value = 0.04100853151850036
$ws.Range("T3").Value = 0.04100853151850037

$ws.Range("I4").Value = 0.08983953209358128
$ws.Range("J4").Value = 0.0898395320935813
$ws.Range("M4").Value = 5.654344666666667
$ws.Range("N4").Value = 16.963034
$ws.Range("O4").Value = 0.2289906587711778
$ws.Range("P4").Value = 0.2289906587711778
$ws.Range("Q4").Value = 0.06774470345133334
$ws.Range("R4").Value = 0.6097023310620001
$ws.Range("S4").Value = 0.02057241363780355
$ws.Range("T4").Value = 0.02057241363780355

$ws.Range("I5").Value = 0.08983953209358128
$ws.Range("J5").Value = 0.0898395320935813
$ws.Range("M5").Value = 0.819389
$ws.Range("N5").Value = 2.458167
$ws.Range("O5").Value = 0.03318376186120772
$ws.Range("P5").Value = 0.03318376186120772
$ws.Range("Q5").Value = 0.009817099609000001
$ws.Range("R5").Value = 0.088353896481
$ws.Range("S5").Value = 0.00298121363871573
$ws.Range("T5").Value = 0.00298121363871573

$ws.Range("I6").Value = 0.08983953209358128
$ws.Range("J6").Value = 0.0898395320935813
$ws.Range("M6").Value = 3.080288333333333
$ws.Range("N6").Value = 9.240864999999999
$ws.Range("O6").Value = 0.1247460662971919
$ws.Range("P6").Value = 0.1247460662971919
$ws.Range("Q6").Value = 0.03690493452166667
$ws.Range("R6").Value = 0.332144410695
$ws.Range("S6").Value = 0.01120712822665459
$ws.Range("T6").Value = 0.01120712822665459

$ws.Range("G7").Value = 0.121379
$ws.Range("H7").Value = 0.364137
$ws.Range("I7").Value = 0.9101604679064187
$ws.Range("J7").Value = 0.9101604679064187
$ws.Range("M7").Value = 3.867218333333334
$ws.Range("N7").Value = 11.601655
$ws.Range("O7").Value = 0.1566152977872902
$ws.Range("P7").Value = 0.1566152977872902
$ws.Range("Q7").Value = 0.4693990940816667
$ws.Range("R7").Value = 4.224591846735001
$ws.Range("S7").Value = 0.1425450527153832
$ws.Range("T7").Value = 0.1425450527153832

$ws.Range("G8").Value = 0.121379
$ws.Range("H8").Value = 0.364137
$ws.Range("I8").Value = 0.9101604679064187
$ws.Range("J8").Value = 0.9101604679064187
$ws.Range("N8").Value = 33.813685
$ws.Range("O8").Value = 0.4564642152831324
$ws.Range("P8").Value = 0.4564642152831324
$ws.Range("Q8").Value = 1.368090423871667
$ws.Range("R8").Value = 12.312813814845
$ws.Range("S8").Value = 0.415455683764632
$ws.Range("T8").Value = 0.415455683764632

$ws.Range("G9").Value = 0.121379
$ws.Range("H9").Value = 0.364137
$ws.Range("I9").Value = 0.9101604679064187
$ws.Range("J9").Value = 0.9101604679064187
$ws.Range("M9").Value = 5.654344666666667
$ws.Range("N9").Value = 16.963034
$ws.Range("O9").Value = 0.2289906587711778
$ws.Range("P9").Value = 0.2289906587711778
$ws.Range("Q9").Value = 0.6863187012953333
$ws.Range("R9").Value = 6.176868311658
$ws.Range("S9").Value = 0.2084182451333743
$ws.Range("T9").Value = 0.2084182451333743

$ws.Range("G10").Value = 0.121379
$ws.Range("H10").Value = 0.364137
$ws.Range("I10").Value = 0.9101604679064187
$ws.Range("J10").Value = 0.9101604679064187
$ws.Range("M10").Value = 0.819389
$ws.Range("N10").Value = 2.458167
$ws.Range("O10").Value = 0.03318376186120772
$ws.Range("P10").Value = 0.03318376186120772
$ws.Range("Q10").Value = 0.09945661743100001
$ws.Range("R10").Value = 0.895109556879
$ws.Range("S10").Value = 0.03020254822249199
$ws.Range("T10").Value = 0.03020254822249199

$ws.Range("G11").Value = 0.121379
$ws.Range("H11").Value = 0.364137
$ws.Range("I11").Value = 0.9101604679064187
$ws.Range("J11").Value = 0.9101604679064187
$ws.Range("M11").Value = 3.080288333333333
$ws.Range("N11").Value = 9.240864999999999
$ws.Range("O11").Value = 0.1247460662971919
$ws.Range("P11").Value = 0.1247460662971919
$ws.Range("Q11").Value = 0.3738823176116667
$ws.Range("R11").Value = 3.364940858505
$ws.Range("S11").Value = 0.1135389380705373
$ws.Range("T11").Value = 0.1135389380705373
